$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("B27").Value = 6183395
$ws.Range("C27").Value = "Finland Veikkausliiga"
$ws.Range("D27").Value = "Finland Veikkausliiga"
$ws.Range("E27").Value = 45108.45833333334
$ws.Range("F27").Value = "VPS Vaasa"
$ws.Range("G27").Value = "FC Inter"
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = "H"
$ws.Range("K27").Value = 3
$ws.Range("L27").Value = 3.1
$ws.Range("M27").Value = 2.2
$ws.Range("N27").Value = 2.75
$ws.Range("O27").Value = 3.1
$ws.Range("P27").Value = 2.375
$ws.Range("Q27").Value = 0
$ws.Range("R27").Value = 2.05
$ws.Range("S27").Value = 1.8
$ws.Range("T27").Value = 2.5
$ws.Range("U27").Value = 1.95
$ws.Range("V27").Value = 1.9
$ws.Range("W27").Value = 1.75
$ws.Range("X27").Value = -1
$ws.Range("Y27").Value = -1
$ws.Range("Z27").Value = 1.05
$ws.Range("AA27").Value = -1
$ws.Range("AB27").Value = 0.95
$ws.Range("AC27").Value = -1

# Row 29
$ws.Range("B29").Value = 6183396
$ws.Range("C29").Value = "Finland Veikkausliiga"
$ws.Range("D29").Value = "Finland Veikkausliiga"
$ws.Range("E29").Value = 45108.45833333334
$ws.Range("F29").Value = "AC Oulu"
$ws.Range("G29").Value = "FC Haka"
$ws.Range("H29").Value = 3
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = "H"
$ws.Range("K29").Value = 2
$ws.Range("L29").Value = 3.1
$ws.Range("M29").Value = 3.4
$ws.Range("N29").Value = 2.3
$ws.Range("O29").Value = 3.1
$ws.Range("P29").Value = 2.75
$ws.Range("Q29").Value = -0.25
$ws.Range("R29").Value = 2.1
$ws.Range("S29").Value = 1.775
$ws.Range("T29").Value = 2.5
$ws.Range("U29").Value = 2
$ws.Range("V29").Value = 1.85
$ws.Range("W29").Value = 1.3
$ws.Range("X29").Value = -1
$ws.Range("Y29").Value = -1
$ws.Range("Z29").Value = 1.1
$ws.Range("AA29").Value = -1
$ws.Range("AB29").Value = 1
$ws.Range("AC29").Value = -1

# Row 56
$ws.Range("B56").Value = 6183421
$ws.Range("C56").Value = "Finland Veikkausliiga"
$ws.Range("D56").Value = "Finland Veikkausliiga"
$ws.Range("E56").Value = 45145.5
$ws.Range("F56").Value = "FC Lahti"
$ws.Range("G56").Value = "FC Honka"
$ws.Range("H56").Value = 1
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = "H"
$ws.Range("K56").Value = 4.5
$ws.Range("L56").Value = 3.5
$ws.Range("M56").Value = 1.8
$ws.Range("N56").Value = 3.75
$ws.Range("O56").Value = 3.3
$ws.Range("P56").Value = 2.05
$ws.Range("Q56").Value = 0.5
$ws.Range("R56").Value = 1.825
$ws.Range("S56").Value = 2.025
$ws.Range("T56").Value = 2.25
$ws.Range("U56").Value = 1.95
$ws.Range("V56").Value = 1.9
$ws.Range("W56").Value = 2.75
$ws.Range("X56").Value = -1
$ws.Range("Y56").Value = -1
$ws.Range("Z56").Value = 0.825
$ws.Range("AA56").Value = -1
$ws.Range("AB56").Value = -1
$ws.Range("AC56").Value = 0.8999999999999999

# Row 57
$ws.Range("B57").Value = 6183420
$ws.Range("C57").Value = "Finland Veikkausliiga"
$ws.Range("D57").Value = "Finland Veikkausliiga"
$ws.Range("E57").Value = 45145.5
$ws.Range("F57").Value = "KTP"
$ws.Range("G57").Value = "FC Inter"
$ws.Range("H57").Value = 1
$ws.Range("I57").Value = 4
$ws.Range("J57").Value = "A"
$ws.Range("K57").Value = 3
$ws.Range("L57").Value = 3.5
$ws.Range("M57").Value = 2.25
$ws.Range("N57").Value = 3.1
$ws.Range("O57").Value = 3.6
$ws.Range("P57").Value = 2.25
$ws.Range("Q57").Value = 0.25
$ws.Range("R57").Value = 1.85
$ws.Range("S57").Value = 2
$ws.Range("T57").Value = 2.75
$ws.Range("U57").Value = 2
$ws.Range("V57").Value = 1.85
$ws.Range("W57").Value = -1
$ws.Range("X57").Value = -1
$ws.Range("Y57").Value = 1.25
$ws.Range("Z57").Value = -1
$ws.Range("AA57").Value = 1
$ws.Range("AB57").Value = 1
$ws.Range("AC57").Value = -1

# Row 83
$ws.Range("B83").Value = 7183904
$ws.Range("C83").Value = "Finland Veikkausliiga"
$ws.Range("D83").Value = "Finland Veikkausliiga"
$ws.Range("E83").Value = 45184.5
$ws.Range("F83").Value = "FC Haka"
$ws.Range("G83").Value = "FC Ilves"
$ws.Range("H83").Value = 2
$ws.Range("I83").Value = 2
$ws.Range("J83").Value = "D"
$ws.Range("K83").Value = 2
$ws.Range("L83").Value = 3.6
$ws.Range("M83").Value = 3.4
$ws.Range("N83").Value = 2.8
$ws.Range("O83").Value = 3.4
$ws.Range("P83").Value = 2.3
$ws.Range("Q83").Value = 0.25
$ws.Range("R83").Value = 1.8
$ws.Range("S83").Value = 2.05
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 2.05
$ws.Range("V83").Value = 1.8
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = 2.4
$ws.Range("Y83").Value = -1
$ws.Range("Z83").Value = 0.4
$ws.Range("AA83").Value = -0.5
$ws.Range("AB83").Value = 1.05
$ws.Range("AC83").Value = -1

# Row 84
$ws.Range("B84").Value = 7183905
$ws.Range("C84").Value = "Finland Veikkausliiga"
$ws.Range("D84").Value = "Finland Veikkausliiga"
$ws.Range("E84").Value = 45184.5
$ws.Range("F84").Value = "FC Lahti"
$ws.Range("G84").Value = "KTP"
$ws.Range("H84").Value = 2
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = "H"
$ws.Range("K84").Value = 2.1
$ws.Range("L84").Value = 3.5
$ws.Range("M84").Value = 3.1
$ws.Range("N84").Value = 2.25
$ws.Range("O84").Value = 3.3
$ws.Range("P84").Value = 3
$ws.Range("Q84").Value = -0.25
$ws.Range("R84").Value = 2
$ws.Range("S84").Value = 1.85
$ws.Range("T84").Value = 2.5
$ws.Range("U84").Value = 2.05
$ws.Range("V84").Value = 1.8
$ws.Range("W84").Value = 1.25
$ws.Range("X84").Value = -1
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 1
$ws.Range("AA84").Value = -1
$ws.Range("AB84").Value = 1.05
$ws.Range("AC84").Value = -1

# Row 90
$ws.Range("B90").Value = 7183919
$ws.Range("C90").Value = "Finland Veikkausliiga"
$ws.Range("D90").Value = "Finland Veikkausliiga"
$ws.Range("E90").Value = 45191.5
$ws.Range("F90").Value = "KTP"
$ws.Range("G90").Value = "FC Ilves"
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 3
$ws.Range("J90").Value = "A"
$ws.Range("K90").Value = 2.625
$ws.Range("L90").Value = 3.4
$ws.Range("M90").Value = 2.55
$ws.Range("N90").Value = 3
$ws.Range("O90").Value = 3.5
$ws.Range("P90").Value = 2.3
$ws.Range("Q90").Value = 0.25
$ws.Range("R90").Value = 1.825
$ws.Range("S90").Value = 2.025
$ws.Range("T90").Value = 2.5
$ws.Range("U90").Value = 1.925
$ws.Range("V90").Value = 1.925
$ws.Range("W90").Value = -1
$ws.Range("X90").Value = -1
$ws.Range("Y90").Value = 1.3
$ws.Range("Z90").Value = -1
$ws.Range("AA90").Value = 1.025
$ws.Range("AB90").Value = 0.925
$ws.Range("AC90").Value = -1

# Row 91
$ws.Range("B91").Value = 7183907
$ws.Range("C91").Value = "Finland Veikkausliiga"
$ws.Range("D91").Value = "Finland Veikkausliiga"
$ws.Range("E91").Value = 45191.5
$ws.Range("F91").Value = "FC Honka"
$ws.Range("G91").Value = "VPS Vaasa"
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 3
$ws.Range("J91").Value = "A"
$ws.Range("K91").Value = 1.909
$ws.Range("L91").Value = 3.5
$ws.Range("M91").Value = 3.9
$ws.Range("N91").Value = 2.25
$ws.Range("O91").Value = 3.3
$ws.Range("P91").Value = 3.25
$ws.Range("Q91").Value = -0.25
$ws.Range("R91").Value = 1.925
$ws.Range("S91").Value = 1.925
$ws.Range("T91").Value = 2.25
$ws.Range("U91").Value = 1.875
$ws.Range("V91").Value = 1.975
$ws.Range("W91").Value = -1
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = 2.25
$ws.Range("Z91").Value = -1
$ws.Range("AA91").Value = 0.925
$ws.Range("AB91").Value = 0.875
$ws.Range("AC91").Value = -1

# Row 110
$ws.Range("B110").Value = 7183917
$ws.Range("C110").Value = "Finland Veikkausliiga"
$ws.Range("D110").Value = "Finland Veikkausliiga"
$ws.Range("E110").Value = 45220.45833333334
$ws.Range("F110").Value = "FC Honka"
$ws.Range("G110").Value = "FC Inter"
$ws.Range("H110").Value = 3
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = "H"
$ws.Range("K110").Value = 2
$ws.Range("L110").Value = 3.5
$ws.Range("M110").Value = 3.2
$ws.Range("N110").Value = 1.6
$ws.Range("O110").Value = 4
$ws.Range("P110").Value = 4.5
$ws.Range("Q110").Value = -0.75
$ws.Range("R110").Value = 1.8
$ws.Range("S110").Value = 2.05
$ws.Range("T110").Value = 2.75
$ws.Range("U110").Value = 1.875
$ws.Range("V110").Value = 1.975
$ws.Range("W110").Value = 0.6000000000000001
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 0.8
$ws.Range("AA110").Value = -1
$ws.Range("AB110").Value = 0.4375
$ws.Range("AC110").Value = -0.5

# Row 111
$ws.Range("B111").Value = 7183918
$ws.Range("C111").Value = "Finland Veikkausliiga"
$ws.Range("D111").Value = "Finland Veikkausliiga"
$ws.Range("E111").Value = 45220.45833333334
$ws.Range("F111").Value = "SJK"
$ws.Range("G111").Value = "VPS Vaasa"
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = "A"
$ws.Range("K111").Value = 2.3
$ws.Range("L111").Value = 3.4
$ws.Range("M111").Value = 2.75
$ws.Range("N111").Value = 2.5
$ws.Range("O111").Value = 3.6
$ws.Range("P111").Value = 2.4
$ws.Range("Q111").Value = 0
$ws.Range("R111").Value = 2.025
$ws.Range("S111").Value = 1.825
$ws.Range("T111").Value = 3
$ws.Range("U111").Value = 1.975
$ws.Range("V111").Value = 1.875
$ws.Range("W111").Value = -1
$ws.Range("X111").Value = -1
$ws.Range("Y111").Value = 1.4
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0.825
$ws.Range("AB111").Value = -1
$ws.Range("AC111").Value = 0.875

# Row 113
$ws.Range("B113").Value = 7380222
$ws.Range("C113").Value = "Finland Veikkausliiga"
$ws.Range("D113").Value = "Finland Veikkausliiga"
$ws.Range("E113").Value = 45224.5
$ws.Range("F113").Value = "FC Honka"
$ws.Range("G113").Value = "FC Inter"
$ws.Range("H113").Value = 4
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = "H"
$ws.Range("K113").Value = 2
$ws.Range("L113").Value = 3.25
$ws.Range("M113").Value = 3.6
$ws.Range("N113").Value = 2.1
$ws.Range("O113").Value = 3.1
$ws.Range("P113").Value = 3.5
$ws.Range("Q113").Value = -0.25
$ws.Range("R113").Value = 1.825
$ws.Range("S113").Value = 2.025
$ws.Range("T113").Value = 2.25
$ws.Range("U113").Value = 1.825
$ws.Range("V113").Value = 2.025
$ws.Range("W113").Value = 1.1
$ws.Range("X113").Value = -1
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = 0.825
$ws.Range("AA113").Value = -1
$ws.Range("AB113").Value = 0.825
$ws.Range("AC113").Value = -1

# Row 114
$ws.Range("B114").Value = 7380223
$ws.Range("C114").Value = "Finland Veikkausliiga"
$ws.Range("D114").Value = "Finland Veikkausliiga"
$ws.Range("E114").Value = 45224.5
$ws.Range("F114").Value = "SJK"
$ws.Range("G114").Value = "AC Oulu"
$ws.Range("H114").Value = 6
$ws.Range("I114").Value = 7
$ws.Range("J114").Value = "A"
$ws.Range("K114").Value = 2.5
$ws.Range("L114").Value = 3.25
$ws.Range("M114").Value = 2.6
$ws.Range("N114").Value = 1.65
$ws.Range("O114").Value = 3.75
$ws.Range("P114").Value = 4
$ws.Range("Q114").Value = -0.75
$ws.Range("R114").Value = 1.95
$ws.Range("S114").Value = 1.9
$ws.Range("T114").Value = 2.5
$ws.Range("U114").Value = 1.85
$ws.Range("V114").Value = 2
$ws.Range("W114").Value = -1
$ws.Range("X114").Value = -1
$ws.Range("Y114").Value = 3
$ws.Range("Z114").Value = -1
$ws.Range("AA114").Value = 0.8999999999999999
$ws.Range("AB114").Value = 0.8500000000000001
$ws.Range("AC114").Value = -1
